$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "# de telefono" header string (was "\# de telefono")
$ws.Cells.Item(1, 2).Value = "# de telefono"

# Add a new value (two spaces) into D1, and move the active selection there
$ws.Range("D1").Value = "  "
$ws.Range("D1").Select()
